# Weekly update: insert this week's newest Perejil price row at the top of
# the data block (row 125) and push the rest of the historical rows down by
# one, matching the "Fruta / hortaliza, semanal" refresh pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 125:159 down to 126:160, carrying formatting along.
$ws.Rows.Item(125).Insert()

# Populate the newly-inserted row 125 with this week's record. All the
# "static" descriptive columns are identical for every row in this sheet.
$ws.Range("A125").Value = 8
$ws.Range("B125").Value = "Terminal La Palmera de La Serena"
$ws.Range("C125").Value = "Coquimbo"
$ws.Range("D125").Value = 44782
$ws.Range("E125").Value = 4
$ws.Range("F125").Value = 100112044
$ws.Range("G125").Value = "Perejil"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 2600
$ws.Range("K125").Value = 2000
$ws.Range("L125").Value = 2500
$ws.Range("M125").Value = 2250
$ws.Range("N125").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O125").Value = "Provincia del Elquí"
$ws.Range("P125").Value = 1500
$ws.Range("Q125").Value = 1.5
$ws.Range("R125").Value = "Hortaliza"
